$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.381.80"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.032.87"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'229.80"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'56.31"
$ws.Range("E8").Value = "  +3.08%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'0.0797"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "2.334.45"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'14.43"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'20.35"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "2.035.67"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "37.356.67"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").Value = "'6.22"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'223.83"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'164.79"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +6.49%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "'2.02"
$ws.Range("E34").Value = "  +9.88%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'5.74"
$ws.Range("E37").Value = "  +8.93%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'3.21"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "1.467.18"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0928"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'94.69"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.82"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.25"
$ws.Range("E45").Value = "  +18.64%  "
$ws.Range("D46").Value = "'16.29"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'7.12"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "2.222.80"
$ws.Range("E51").Value = "  +0.48%  "
